$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 15738.63
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 16328.577
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 48985.731
$ws.Range("M46").Value = -1081
$ws.Range("N46").Value = -49223.731

$ws.Range("H60").Value = 15738.63
$ws.Range("I60").Value = 400
$ws.Range("J60").Value = 16328.577
$ws.Range("K60").Value = 1200
$ws.Range("L60").Value = 48985.731
$ws.Range("M60").Value = -716
$ws.Range("N60").Value = -49953.731

$ws.Range("H64").Value = 4070.3064
$ws.Range("I64").Value = 3429.2683
$ws.Range("J64").Value = 5321.857
$ws.Range("K64").Value = 3429.2683
$ws.Range("L64").Value = 5321.857
$ws.Range("M64").Value = -3181.2683
$ws.Range("N64").Value = -5817.857

$ws.Range("H67").Value = 4070.3064
$ws.Range("I67").Value = 3429.2683
$ws.Range("J67").Value = 5321.857
$ws.Range("K67").Value = 3429.2683
$ws.Range("L67").Value = 5321.857
$ws.Range("M67").Value = -2571.2683
$ws.Range("N67").Value = -7037.857

$ws.Range("H86").Value = 200745
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 200745
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 200745
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -202991

$ws.Range("H89").Value = 200745
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 200745
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 1003725
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -1014957

$ws.Range("H98").Value = 1312
$ws.Range("I98").Value = 1297.7778
$ws.Range("J98").Value = 1333.3334
$ws.Range("K98").Value = 1297.7778
$ws.Range("L98").Value = 1333.3334
$ws.Range("M98").Value = 200.2221999999999
$ws.Range("N98").Value = -4329.3334

$ws.Range("H122").Value = 1312
$ws.Range("I122").Value = 1297.7778
$ws.Range("J122").Value = 1333.3334
$ws.Range("K122").Value = 3893.3334
$ws.Range("L122").Value = 4000.0002
$ws.Range("M122").Value = -1443.3334
$ws.Range("N122").Value = -8900.0002

$ws.Range("H125").Value = 2746.25
$ws.Range("I125").Value = 2945
$ws.Range("J125").Value = 2547.5
$ws.Range("K125").Value = 26505
$ws.Range("L125").Value = 22927.5
$ws.Range("M125").Value = -24045
$ws.Range("N125").Value = -27847.5

$ws.Range("H129").Value = 876.19354
$ws.Range("J129").Value = 891.98303
$ws.Range("L129").Value = 2675.94909
$ws.Range("N129").Value = -12675.94909

$ws.Range("H137").Value = 2664.4285
$ws.Range("I137").Value = 1999.8334
$ws.Range("J137").Value = 3162.875
$ws.Range("K137").Value = 5999.5002
$ws.Range("L137").Value = 9488.625
$ws.Range("M137").Value = -3449.5002
$ws.Range("N137").Value = -14588.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3812.0513
$ws.Range("I32").Value = 2372.2112
$ws.Range("J32").Value = 18416.143
$ws.Range("K32").Value = 2372.2112
$ws.Range("L32").Value = 18416.143
$ws.Range("M32").Value = -2085.2112
$ws.Range("N32").Value = -18990.143

$ws.Range("H61").Value = 20014
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 20014
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 20014
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -20438

$ws.Range("H74").Value = 1176.52
$ws.Range("I74").Value = 1116
$ws.Range("J74").Value = 1267.3
$ws.Range("K74").Value = 1116
$ws.Range("L74").Value = 1267.3
$ws.Range("M74").Value = -242
$ws.Range("N74").Value = -3015.3

$ws.Range("H77").Value = 1176.52
$ws.Range("I77").Value = 1116
$ws.Range("J77").Value = 1267.3
$ws.Range("K77").Value = 5580
$ws.Range("L77").Value = 6336.5
$ws.Range("M77").Value = -1212
$ws.Range("N77").Value = -15072.5

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("N109").ClearContents()

$ws.Range("H112").Value = 7942.2
$ws.Range("J112").Value = 7942.2
$ws.Range("L112").Value = 7942.2
$ws.Range("N112").Value = -10896.2

$ws.Range("H114").Value = 36099.5
$ws.Range("J114").Value = 36099.5
$ws.Range("L114").Value = 36099.5
$ws.Range("N114").Value = -44777.5

$ws.Range("H132").Value = 2764.147
$ws.Range("I132").Value = 1207.4706
$ws.Range("J132").Value = 4320.8237
$ws.Range("K132").Value = 3622.4118
$ws.Range("L132").Value = 12962.4711
$ws.Range("M132").Value = -1092.4118
$ws.Range("N132").Value = -18022.4711

$ws.Range("H136").Value = 20014
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 20014
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 60042
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -65142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2342.5
$ws.Range("I86").Value = 2744
$ws.Range("K86").Value = 2744
$ws.Range("M86").Value = -1621

$ws.Range("H89").Value = 2342.5
$ws.Range("I89").Value = 2744
$ws.Range("K89").Value = 13720
$ws.Range("M89").Value = -8104

$ws.Range("H134").Value = 2671.2666
$ws.Range("I134").Value = 1188.125
$ws.Range("J134").Value = 4366.2856
$ws.Range("K134").Value = 3564.375
$ws.Range("L134").Value = 13098.8568
$ws.Range("M134").Value = -1029.375
$ws.Range("N134").Value = -18168.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2227.8
$ws.Range("I31").Value = 1992.5385
$ws.Range("J31").Value = 3757
$ws.Range("K31").Value = 1992.5385
$ws.Range("L31").Value = 3757
$ws.Range("M31").Value = -1697.5385
$ws.Range("N31").Value = -4347

$ws.Range("H34").Value = 2227.8
$ws.Range("I34").Value = 1992.5385
$ws.Range("J34").Value = 3757
$ws.Range("K34").Value = 1992.5385
$ws.Range("L34").Value = 3757
$ws.Range("M34").Value = -1790.5385
$ws.Range("N34").Value = -4161

$ws.Range("H58").Value = 2240.0303
$ws.Range("I58").Value = 1524.7142
$ws.Range("J58").Value = 2767.1052
$ws.Range("K58").Value = 1524.7142
$ws.Range("L58").Value = 2767.1052
$ws.Range("M58").Value = -1321.7142
$ws.Range("N58").Value = -3173.1052

$ws.Range("H136").Value = 2240.0303
$ws.Range("I136").Value = 1524.7142
$ws.Range("J136").Value = 2767.1052
$ws.Range("K136").Value = 4574.142599999999
$ws.Range("L136").Value = 8301.3156
$ws.Range("M136").Value = -2024.142599999999
$ws.Range("N136").Value = -13401.3156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5785.2
$ws.Range("I3").Value = 6020
$ws.Range("J3").Value = 5433
$ws.Range("K3").Value = 18060
$ws.Range("L3").Value = 16299
$ws.Range("M3").Value = -17948
$ws.Range("N3").Value = -16523

$ws.Range("H132").Value = 1001.64703
$ws.Range("I132").Value = 744.3333
$ws.Range("J132").Value = 1291.125
$ws.Range("K132").Value = 6698.9997
$ws.Range("L132").Value = 11620.125
$ws.Range("M132").Value = -4168.9997
$ws.Range("N132").Value = -16680.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H80").Value = 1732.5
$ws.Range("I80").Value = 2062.2
$ws.Range("J80").Value = 1549.3334
$ws.Range("K80").Value = 2062.2
$ws.Range("L80").Value = 1549.3334
$ws.Range("M80").Value = -1064.2
$ws.Range("N80").Value = -3545.3334

$ws.Range("H83").Value = 1732.5
$ws.Range("I83").Value = 2062.2
$ws.Range("J83").Value = 1549.3334
$ws.Range("K83").Value = 10311
$ws.Range("L83").Value = 7746.666999999999
$ws.Range("M83").Value = -5319
$ws.Range("N83").Value = -17730.667

$ws.Range("H132").Value = 5347
$ws.Range("J132").Value = 7702
$ws.Range("L132").Value = 23106
$ws.Range("N132").Value = -28166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17094.766
$ws.Range("I132").Value = 32229.715
$ws.Range("J132").Value = 6500.3
$ws.Range("K132").Value = 96689.145
$ws.Range("L132").Value = 19500.9
$ws.Range("M132").Value = -94159.145
$ws.Range("N132").Value = -24560.9

$ws.Range("H136").Value = 55563390
$ws.Range("I136").Value = 10499.5
$ws.Range("J136").Value = 83339840
$ws.Range("K136").Value = 31498.5
$ws.Range("L136").Value = 250019520
$ws.Range("M136").Value = -28948.5
$ws.Range("N136").Value = -250024620

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49333.223
$ws.Range("J46").Value = 49333.223
$ws.Range("L46").Value = 49333.223
$ws.Range("N46").Value = -49795.223

$ws.Range("H56").Value = 6057
$ws.Range("J56").Value = 6742.6665
$ws.Range("L56").Value = 6742.6665
$ws.Range("N56").Value = -8170.6665

$ws.Range("H62").Value = 4315.8335
$ws.Range("I62").Value = 2965
$ws.Range("K62").Value = 2965
$ws.Range("M62").Value = -2341

$ws.Range("H65").Value = 4315.8335
$ws.Range("I65").Value = 2965
$ws.Range("K65").Value = 14825
$ws.Range("M65").Value = -11705

$ws.Range("H75").Value = 33000
$ws.Range("J75").Value = 33000
$ws.Range("L75").Value = 33000
$ws.Range("N75").Value = -34872

$ws.Range("H78").Value = 33000
$ws.Range("J78").Value = 33000
$ws.Range("L78").Value = 99000
$ws.Range("N78").Value = -108360

$ws.Range("H132").Value = 2161.7742
$ws.Range("I132").Value = 1696.8077
$ws.Range("J132").Value = 4579.6
$ws.Range("K132").Value = 5090.4231
$ws.Range("L132").Value = 13738.8
$ws.Range("M132").Value = -2560.4231
$ws.Range("N132").Value = -18798.8

$ws.Range("H134").Value = 49333.223
$ws.Range("J134").Value = 49333.223
$ws.Range("L134").Value = 147999.669
$ws.Range("N134").Value = -153069.669

$ws.Range("H136").Value = 9225.632
$ws.Range("J136").Value = 14983.182
$ws.Range("L136").Value = 44949.546
$ws.Range("N136").Value = -50049.546

$ws.Range("H137").Value = 59283.332
$ws.Range("J137").Value = 59283.332
$ws.Range("L137").Value = 59283.332
$ws.Range("N137").Value = -69483.33199999999
